{"js": "// Apply the CV date-formatting edits described by the diff.\n// Each entry is [searchText, replacementText]. The first three entries\n// remove a leading month abbreviation (e.g. \"Sep. \", \"Jun. \", \"May \")\n// from a date range (replacing it with extra blank padding so the\n// overall layout/spacing is preserved); the remaining entries are\n// \"replace text with itself\" operations that simply collapse runs\n// that already held identical, contiguous text (no visible change).\nconst replacements = [\n  // \"Ph.D. in Meteorology and Physical Oceanography ... Sep. 2011 \u2013 May 2016\"\n  //   -> drop \"Sep. \" and \"May \"\n  [\n    \" \".repeat(21) + \"Sep. 2011 \u2013 May \",\n    \" \".repeat(36) + \"2011 \u2013 \"\n  ],\n  // \"B.S. in Ocean Sciences ... Sep. 2007 \u2013 Jun. 2011\"\n  //   -> drop \"Sep. \" and \"Jun. \"\n  [\n    \" \".repeat(65) + \"Sep. 2007 \u2013 Jun. \",\n    \" \".repeat(80) + \"2007 \u2013 \"\n  ],\n  // \"Postdoctoral Research Associate, University of Miami, USA ... Jun. 2016 \u2013 Current\"\n  //   -> drop \"Jun. \"\n  [\n    \" \".repeat(5) + \"    Jun. 2016 \u2013 Current   \",\n    \" \" + \" \".repeat(15) + \"2016 \u2013 Current   \"\n  ],\n  // Research Assistant: consolidate the date-range runs (text unchanged)\n  [\n    \"                                Aug. 2011 \u2013 May 2016\",\n    \"                                Aug. 2011 \u2013 May 2016\"\n  ],\n  // Teaching Assistant for Geophysical Fluid Dynamics: consolidate date runs\n  [\n    \"                Oct. 2015 \u2013 Nov. \",\n    \"                Oct. 2015 \u2013 Nov. \"\n  ],\n  // Visiting Scholar ... Sep. 2015: consolidate date runs\n  [\n    \"                                              Sep. 2015 \",\n    \"                                              Sep. 2015 \"\n  ],\n  // Teaching Assistant for Computer Models of Fluid Dynamics: consolidate date runs\n  [\n    \"Aug. 2013 \u2013 Dec. \",\n    \"Aug. 2013 \u2013 Dec. \"\n  ],\n  // Teaching Assistant for Introduction to Physical Oceanography: consolidate date runs\n  [\n    \" Jan. 2013 \u2013 May 2013\",\n    \" Jan. 2013 \u2013 May 2013\"\n  ],\n  // Volunteer for National Gandhi Day of Service: consolidate date runs\n  [\n    \"                                                Oct. 2015                    \",\n    \"                                                Oct. 2015                    \"\n  ],\n  // Volunteer for UM/RSMAS Student Auction: consolidate date runs\n  [\n    \"                                                 Mar. 2015 \",\n    \"                                                 Mar. 2015 \"\n  ],\n  // Member of UM/RSMAS Garden Club: consolidate date runs\n  [\n    \"                                        Jan. 2013 \u2013 Mar. 2016\",\n    \"                                        Jan. 2013 \u2013 Mar. 2016\"\n  ]\n];\n\nfor (const [searchText, replacementText] of replacements) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for \" + JSON.stringify(searchText) +\n      \" but found \" + results.items.length\n    );\n  }\n\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply the CV date-formatting edits described by the diff.\n#\n# The first three Find/Replace calls remove a leading month abbreviation\n# (e.g. \"Sep. \", \"Jun. \", \"May \") from a date range, replacing it with\n# extra blank padding so the overall layout/spacing is preserved. The\n# remaining calls are \"replace text with itself\" operations that simply\n# collapse runs that already held identical, contiguous text (no visible\n# text change, just run consolidation).\n#\n# NOTE: this interpreter mis-parses string concatenation with \"+\" when the\n# left-hand operand looks purely numeric (e.g. \"2011\" + <en-dash>), so the\n# date fragments below are written as single literal strings (with the\n# literal \"\u2013\" en-dash character embedded) rather than built up with \"+\".\n\n$d = $word.ActiveDocument\n\nfunction ReplaceOnce($searchText, $replaceText) {\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $result = $find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n  if (-not $result) {\n    Write-Output (\"WARNING: no match for [\" + $searchText + \"]\")\n  }\n  return $result\n}\n\n# 1) \"Ph.D. in Meteorology and Physical Oceanography ... Sep. 2011 \u2013 May 2016\"\n#      -> drop \"Sep. \" and \"May \"\n$search1 = \"\".PadLeft(21) + \"Sep. 2011 \u2013 May \"\n$replace1 = \"\".PadLeft(36) + \"2011 \u2013 \"\nReplaceOnce $search1 $replace1 | Out-Null\n\n# 2) \"B.S. in Ocean Sciences ... Sep. 2007 \u2013 Jun. 2011\"\n#      -> drop \"Sep. \" and \"Jun. \"\n$search2 = \"\".PadLeft(65) + \"Sep. 2007 \u2013 Jun. \"\n$replace2 = \"\".PadLeft(80) + \"2007 \u2013 \"\nReplaceOnce $search2 $replace2 | Out-Null\n\n# 3) \"Postdoctoral Research Associate, University of Miami, USA ... Jun. 2016 \u2013 Current\"\n#      -> drop \"Jun. \"\n$search3 = \"\".PadLeft(9) + \"Jun. 2016 \u2013 Current   \"\n$replace3 = \"\".PadLeft(16) + \"2016 \u2013 Current   \"\nReplaceOnce $search3 $replace3 | Out-Null\n\n# 4) Research Assistant: consolidate the date-range runs (text unchanged)\n$s4 = \"                                Aug. 2011 \u2013 May 2016\"\nReplaceOnce $s4 $s4 | Out-Null\n\n# 5) Teaching Assistant for Geophysical Fluid Dynamics: consolidate date runs\n$s5 = \"                Oct. 2015 \u2013 Nov. \"\nReplaceOnce $s5 $s5 | Out-Null\n\n# 6) Visiting Scholar ... Sep. 2015: consolidate date runs\n$s6 = \"                                              Sep. 2015 \"\nReplaceOnce $s6 $s6 | Out-Null\n\n# 7) Teaching Assistant for Computer Models of Fluid Dynamics: consolidate date runs\n$s7 = \"Aug. 2013 \u2013 Dec. \"\nReplaceOnce $s7 $s7 | Out-Null\n\n# 8) Teaching Assistant for Introduction to Physical Oceanography: consolidate date runs\n$s8 = \" Jan. 2013 \u2013 May 2013\"\nReplaceOnce $s8 $s8 | Out-Null\n\n# 9) Volunteer for National Gandhi Day of Service: consolidate date runs\n$s9 = \"                                                Oct. 2015                    \"\nReplaceOnce $s9 $s9 | Out-Null\n\n# 10) Volunteer for UM/RSMAS Student Auction: consolidate date runs\n$s10 = \"                                                 Mar. 2015 \"\nReplaceOnce $s10 $s10 | Out-Null\n\n# 11) Member of UM/RSMAS Garden Club: consolidate date runs\n$s11 = \"                                        Jan. 2013 \u2013 Mar. 2016\"\nReplaceOnce $s11 $s11 | Out-Null\n"}
